$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Determine last used row on the sheet (data starts at row 2)
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

# Column C holds the "Förändrad" (last changed) date serial value.
# Every data row currently stores 45178 (2023-09-09); bump it by one day
# to 45179 (2023-09-10) for all rows.
$range = $ws.Range("C2:C$lastRow")
$range.Value = 45179
